# Applies the cryptos.xlsx data-refresh edit described by the commit diff.
# Column D/E (and occasionally B/C) values are refreshed per-row; all cells
# are plain text (t="inlineStr" in the source), so numeric-looking price
# strings are written with a leading quote-prefix to keep Excel from
# re-typing them as numbers, then the style is reset back to Normal so no
# stray quote-prefix formatting is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.827.18"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.498.03"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'539.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").Value = "'143.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "2.519.06"
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "'5.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.13%  "
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "2.937.07"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").Value = "'23.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "58.759.21"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Value = "2.513.62"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "'323.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "'61.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("E25").Value = "  -4.58%  "
$ws.Range("D26").Value = "'0.162"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("E27").Value = "  +2.31%  "
$ws.Range("D28").Value = "2.615.71"
$ws.Range("E28").Value = "  +2.67%  "
$ws.Range("D29").Value = "'7.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").Value = "0.0₃0768"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").Value = "'1.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.04%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'156.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("E36").Value = "  +3.65%  "
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("D38").Value = "'4.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.63%  "
$ws.Range("D39").Value = "'1.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.01%  "
$ws.Range("D40").Value = "'5.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").Value = "'36.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "'296.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.04%  "
$ws.Range("D43").Value = "'3.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'0.821"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E46").Value = "  +3.97%  "
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.0928"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'123.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.94%  "
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").Value = "'0.0228"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
